$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "stim/" folder prefix from the image filenames in C2:C9
for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 -replace "^stim/", ""
}

# Rename header cell C1 from "image" to "imageFile"
$ws.Range("C1").Value2 = "imageFile"

# Move the active selection to C1 (it was previously on C9)
$ws.Range("C1").Select()
